# Scheduled runner update: refresh market-price / profit columns (H:N)
# for a batch of leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1125.9131
$ws.Range("I98").Value = 1151.909
$ws.Range("J98").Value = 554
$ws.Range("K98").Value = 1151.909
$ws.Range("L98").Value = 554
$ws.Range("M98").Value = 346.0909999999999
$ws.Range("N98").Value = -3550

$ws.Range("H113").Value = 3974.5
$ws.Range("I113").Value = 3674.25
$ws.Range("J113").Value = 4124.625
$ws.Range("K113").Value = 3674.25
$ws.Range("L113").Value = 4124.625
$ws.Range("M113").Value = -420.25
$ws.Range("N113").Value = -10632.625

$ws.Range("H122").Value = 1125.9131
$ws.Range("I122").Value = 1151.909
$ws.Range("J122").Value = 554
$ws.Range("K122").Value = 3455.727
$ws.Range("L122").Value = 1662
$ws.Range("M122").Value = -1005.727
$ws.Range("N122").Value = -6562

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4960.268
$ws.Range("I61").Value = 4975.163
$ws.Range("K61").Value = 4975.163
$ws.Range("M61").Value = -4763.163

$ws.Range("H102").Value = 4131.7334
$ws.Range("J102").Value = 4895.4
$ws.Range("L102").Value = 4895.4
$ws.Range("N102").Value = -8139.4

$ws.Range("H136").Value = 4960.268
$ws.Range("I136").Value = 4975.163
$ws.Range("K136").Value = 14925.489
$ws.Range("M136").Value = -12375.489

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1681.6666
$ws.Range("I86").Value = 1431.2
$ws.Range("K86").Value = 1431.2
$ws.Range("M86").Value = -308.2

$ws.Range("H89").Value = 1681.6666
$ws.Range("I89").Value = 1431.2
$ws.Range("K89").Value = 7156
$ws.Range("M89").Value = -1540

$ws.Range("H134").Value = 3102.7942
$ws.Range("I134").Value = 3031.0688
$ws.Range("K134").Value = 9093.206399999999
$ws.Range("M134").Value = -6558.206399999999

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1984.9062
$ws.Range("I31").Value = 1993.7307
$ws.Range("J31").Value = 1946.6666
$ws.Range("K31").Value = 1993.7307
$ws.Range("L31").Value = 1946.6666
$ws.Range("M31").Value = -1698.7307
$ws.Range("N31").Value = -2536.6666

$ws.Range("H34").Value = 1984.9062
$ws.Range("I34").Value = 1993.7307
$ws.Range("J34").Value = 1946.6666
$ws.Range("K34").Value = 1993.7307
$ws.Range("L34").Value = 1946.6666
$ws.Range("M34").Value = -1791.7307
$ws.Range("N34").Value = -2350.6666

$ws.Range("H58").Value = 2002.2858
$ws.Range("I58").Value = 1570.1333
$ws.Range("K58").Value = 1570.1333
$ws.Range("M58").Value = -1367.1333

$ws.Range("H94").Value = 15359.5
$ws.Range("I94").Value = 28412
$ws.Range("J94").Value = 2307
$ws.Range("K94").Value = 28412
$ws.Range("L94").Value = 2307
$ws.Range("M94").Value = -27961
$ws.Range("N94").Value = -3209

$ws.Range("H99").Value = 10712.366
$ws.Range("I99").Value = 6576.6665
$ws.Range("K99").Value = 6576.6665
$ws.Range("M99").Value = -5078.6665

$ws.Range("H126").Value = 10712.366
$ws.Range("I126").Value = 6576.6665
$ws.Range("K126").Value = 19729.9995
$ws.Range("M126").Value = -17259.9995

$ws.Range("H134").Value = 2080.45
$ws.Range("I134").Value = 2000.8334
$ws.Range("K134").Value = 6002.5002
$ws.Range("M134").Value = -3467.5002

$ws.Range("H136").Value = 2002.2858
$ws.Range("I136").Value = 1570.1333
$ws.Range("K136").Value = 4710.3999
$ws.Range("M136").Value = -2160.3999

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 10326
$ws.Range("I88").Value = 13000
$ws.Range("J88").Value = 8320.5
$ws.Range("K88").Value = 39000
$ws.Range("L88").Value = 24961.5
$ws.Range("M88").Value = -38572
$ws.Range("N88").Value = -25817.5

$ws.Range("H91").Value = 10326
$ws.Range("I91").Value = 13000
$ws.Range("J91").Value = 8320.5
$ws.Range("K91").Value = 39000
$ws.Range("L91").Value = 24961.5
$ws.Range("M91").Value = -37518
$ws.Range("N91").Value = -27925.5

$ws.Range("H104").Value = 2879.75
$ws.Range("I104").Value = 506.33334
$ws.Range("J104").Value = 10000
$ws.Range("K104").Value = 1519.00002
$ws.Range("L104").Value = 30000
$ws.Range("M104").Value = 1101.99998
$ws.Range("N104").Value = -35242

$ws.Range("H109").Value = 1510.35
$ws.Range("I109").Value = 805.2143
$ws.Range("K109").Value = 2415.6429
$ws.Range("M109").Value = -1375.6429

$ws.Range("H131").Value = 685819.4
$ws.Range("I131").Value = 3677746.2
$ws.Range("K131").Value = 11033238.6
$ws.Range("M131").Value = -11028198.6

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 32829
$ws.Range("I46").Value = 2743.5
$ws.Range("J46").Value = 93000
$ws.Range("K46").Value = 2743.5
$ws.Range("L46").Value = 93000
$ws.Range("M46").Value = -2587.5
$ws.Range("N46").Value = -93312

$ws.Range("H70").Value = 7974
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 7974
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H102").Value = 6329.1904
$ws.Range("I102").Value = 6594.2666
$ws.Range("J102").Value = 5666.5
$ws.Range("K102").Value = 6594.2666
$ws.Range("L102").Value = 5666.5
$ws.Range("M102").Value = -4972.2666
$ws.Range("N102").Value = -8910.5

$ws.Range("H122").Value = 3113.963
$ws.Range("I122").Value = 2767.158
$ws.Range("J122").Value = 3937.625
$ws.Range("K122").Value = 8301.474
$ws.Range("L122").Value = 11812.875
$ws.Range("M122").Value = -5851.474
$ws.Range("N122").Value = -16712.875

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3199.6667
$ws.Range("I7").Value = 3199.6667
$ws.Range("K7").Value = 3199.6667
$ws.Range("M7").Value = -3087.6667

$ws.Range("H55").Value = 214.15
$ws.Range("J55").Value = 401.57144
$ws.Range("L55").Value = 401.57144
$ws.Range("N55").Value = -747.5714399999999

$ws.Range("H126").Value = 3199.6667
$ws.Range("I126").Value = 3199.6667
$ws.Range("K126").Value = 9599.000100000001
$ws.Range("M126").Value = -7129.000100000001

$ws.Range("H132").Value = 2920
$ws.Range("I132").Value = 2920
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8760
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6230
$ws.Range("N132").ClearContents()

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2490.9792
$ws.Range("I132").Value = 2103.0256
$ws.Range("J132").Value = 4172.1113
$ws.Range("K132").Value = 6309.0768
$ws.Range("L132").Value = 12516.3339
$ws.Range("M132").Value = -3779.0768
$ws.Range("N132").Value = -17576.3339

$ws.Range("H136").Value = 1449.2391
$ws.Range("I136").Value = 1380.561
$ws.Range("J136").Value = 2012.4
$ws.Range("K136").Value = 4141.683
$ws.Range("L136").Value = 6037.200000000001
$ws.Range("M136").Value = -1591.683
$ws.Range("N136").Value = -11137.2
